$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Formula = '="abc"'
$ws.Range("B6").Formula = "=4711"
$ws.Range("C7").Formula = "=TRUE"
$ws.Range("D8").Formula = "=C7"

$ws.Range("D9").Select()
